# Update sheet "汽車" (car, sheet index 2): replace the row-1 header labels
# with proper column names, and append the extra metadata columns (H:N)
# to both the header row and the data row - matching the layout already
# used on the other property sheets ("建物"/"股票").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 1 - header labels for the existing columns B:G
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "capacity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "register_date"
$ws.Cells.Item(1,6).Value = "register_reason"
$ws.Cells.Item(1,7).Value = "acquire_value"

# Row 1 - new header labels for the appended columns H:N (same bold/border
# style as the rest of the header row)
$ws.Cells.Item(1,8).Value  = "property_category"
$ws.Cells.Item(1,9).Value  = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = "legislator_id"
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = "index"

$headerRow = $ws.Range("H1:N1")
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160
$headerRow.Borders.LineStyle = 1

# Row 2 - new metadata values appended in columns H:N
$ws.Cells.Item(2,8).Value  = "land"
$ws.Cells.Item(2,9).Value  = "normal"
# Force the register/report date to stay plain text instead of being
# auto-converted to a date serial number.
$ws.Cells.Item(2,10).Value = "'2013-11-12"
$ws.Cells.Item(2,10).Style = "Normal"
$ws.Cells.Item(2,11).Value = "王育敏"
$ws.Cells.Item(2,12).Value = 1728
$ws.Cells.Item(2,13).Value = "tmped871"
$ws.Cells.Item(2,14).Value = 29
